$d = $word.ActiveDocument

# 1) Programa - Portuguese paragraph: split merged text into 5 lines separated by manual line breaks
$d.Content.Find.Execute(
  "Perspectivas: nanociência e nanotecnologia - a distinção; Implicações sociais de nanoNanotools: métodos de caracterização; Métodos de fabricaçãoFísica: Propriedades e fenômenos: materiais, estrutura e nanosurface; Energia na nanoescalaQuímica: síntese e modificação: nanomateriais à base de carbono; Interações químicas na nanoescalaAplicações: nanoetronics; nanomagnetismo; nanomecânica",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Perspectivas: nanociência e nanotecnologia - a distinção; Implicações sociais de nano^lNanotools: métodos de caracterização; Métodos de fabricação^lFísica: Propriedades e fenômenos: materiais, estrutura e nanosurface; Energia na nanoescala^lQuímica: síntese e modificação: nanomateriais à base de carbono; Interações químicas na nanoescala^lAplicações: nanoetronics; nanomagnetismo; nanomecânica",
  2
) | Out-Null

# 2) Programa - English paragraph: split merged text into 5 lines separated by manual line breaks
$d.Content.Find.Execute(
  "Perspectives: Nanoscience and Nanotechnology—The Distinction; Societal Implications of NanoNanotools: Characterization Methods; Fabrication MethodsPhysics: Properties and Phenomena: Materials, Structure, and the Nanosurface; Energy at the NanoscaleChemistry: Synthesis and Modification: Carbon-Based Nanomaterials; Chemical Interactions at the NanoscaleApplications: nanoeletronics; nanomagnetism; nanomechanics",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Perspectives: Nanoscience and Nanotechnology—The Distinction; Societal Implications of Nano^lNanotools: Characterization Methods; Fabrication Methods^lPhysics: Properties and Phenomena: Materials, Structure, and the Nanosurface; Energy at the Nanoscale^lChemistry: Synthesis and Modification: Carbon-Based Nanomaterials; Chemical Interactions at the Nanoscale^lApplications: nanoeletronics; nanomagnetism; nanomechanics",
  2
) | Out-Null

# 3) Avaliação - Método run text: split merged text into 4 lines separated by manual line breaks
#    (the trailing <w:br/> already present in the run stays untouched)
$d.Content.Find.Execute(
  "Aulas expositivas e seminários.CritérioDuas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3Norma de RecuperaçãoAplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Aulas expositivas e seminários.^lCritério^lDuas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3^lNorma de Recuperação^lAplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação",
  2
) | Out-Null

# 4) Bibliografia paragraph: split merged text into 3 lines separated by manual line breaks
$d.Content.Find.Execute(
  "Gabor L. Hornyak, H.F. Tibbals, Joydeep Dutta, John J. Moore. Introduction to Nanoscience and Nanotechnology. CRC Press. 2009TIMP, G. Nanotechnology, Springer, 1998.Bhushan, B. (ed.) Springer Handbook of Nanotechnology, Springer, 2010.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Gabor L. Hornyak, H.F. Tibbals, Joydeep Dutta, John J. Moore. Introduction to Nanoscience and Nanotechnology. CRC Press. 2009^lTIMP, G. Nanotechnology, Springer, 1998.^lBhushan, B. (ed.) Springer Handbook of Nanotechnology, Springer, 2010.",
  2
) | Out-Null
